$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EW")

$ws.Range("B4").Value = 768000000.0
$ws.Range("C4").Value = 802000000.0
$ws.Range("D4").Value = 773000000.0
$ws.Range("E4").Value = 735000000.0
$ws.Range("F4").Value = 662000000.0

$ws.Range("B15").Value = 150000000.0
$ws.Range("C15").Value = 197000000.0
$ws.Range("D15").Value = 163000000.0
$ws.Range("E15").Value = 153000000.0
$ws.Range("F15").Value = 152000000.0

$ws.Range("B22").Value = -209000000.0
$ws.Range("C22").Value = -231000000.0
$ws.Range("D22").Value = -215000000.0
$ws.Range("E22").Value = -204000000.0
$ws.Range("F22").Value = -158000000.0
